$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.032261480122114
$ws.Range("D2").Value = 1.035308803580233
$ws.Range("E2").Value = 0.9926147277508489
$ws.Range("F2").Value = 1.041942261124692
$ws.Range("I2").Value = 1.036631949492165
$ws.Range("J2").Value = 1.037392672141993
$ws.Range("K2").Value = 1.038105779121182
$ws.Range("L2").Value = 0.9955398523336033
$ws.Range("M2").Value = 1.044720349950321
$ws.Range("N2").Value = 1.016427378699448

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.033182000603881
$ws.Range("D3").Value = 1.035984231674604
$ws.Range("E3").Value = 0.9936372048519304
$ws.Range("F3").Value = 1.043105570518958
$ws.Range("I3").Value = 1.036863523839338
$ws.Range("J3").Value = 1.037955560540912
$ws.Range("K3").Value = 1.038590931600067
$ws.Range("L3").Value = 0.9963617723202692
$ws.Range("M3").Value = 1.045693465348894
$ws.Range("N3").Value = 1.016615980641403

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.033777730515216
$ws.Range("D4").Value = 1.036421361334773
$ws.Range("E4").Value = 0.9942998659930995
$ws.Range("F4").Value = 1.043858823678255
$ws.Range("I4").Value = 1.037012260446804
$ws.Range("J4").Value = 1.038319259017732
$ws.Range("K4").Value = 1.038904260919876
$ws.Range("L4").Value = 0.9968940712668345
$ws.Range("M4").Value = 1.04632307533639
$ws.Range("N4").Value = 1.016737787491219

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.034028196712646
$ws.Range("D5").Value = 1.036605149372438
$ws.Range("E5").Value = 0.9945786998346017
$ws.Range("F5").Value = 1.044175614520164
$ws.Range("I5").Value = 1.037074523849884
$ws.Range("J5").Value = 1.038472030883515
$ws.Range("K5").Value = 1.039035841068437
$ws.Range("L5").Value = 0.997117960005301
$ws.Range("M5").Value = 1.046587748817693
$ws.Range("N5").Value = 1.016788939485274

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.034070252323809
$ws.Range("D6").Value = 1.036636009261034
$ws.Range("E6").Value = 0.9946255319796338
$ws.Range("F6").Value = 1.044228812325243
$ws.Range("I6").Value = 1.037084962572465
$ws.Range("J6").Value = 1.038497674489349
$ws.Range("K6").Value = 1.039057925521332
$ws.Range("L6").Value = 0.9971555583673453
$ws.Range("M6").Value = 1.046632187805837
$ws.Range("N6").Value = 1.016797524862097

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.03378107717367
$ws.Range("D7").Value = 1.036423817047229
$ws.Range("E7").Value = 0.9943035907982488
$ws.Range("F7").Value = 1.043863056168687
$ws.Range("I7").Value = 1.037013093456959
$ws.Range("J7").Value = 1.038321300862524
$ws.Range("K7").Value = 1.038906019664895
$ws.Range("L7").Value = 0.9968970624462087
$ws.Range("M7").Value = 1.046326611973612
$ws.Range("N7").Value = 1.016738471205475

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.032572555007082
$ws.Range("D8").Value = 1.035537049964492
$ws.Range("E8").Value = 0.9929600610674301
$ws.Range("F8").Value = 1.042335300581817
$ws.Range("I8").Value = 1.036710440237049
$ws.Range("J8").Value = 1.037583011879963
$ws.Range("K8").Value = 1.038269861861485
$ws.Range("L8").Value = 0.995817528259106
$ws.Range("M8").Value = 1.045049231423173
$ws.Range("N8").Value = 1.016491165332899

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.03044371182344
$ws.Range("D9").Value = 1.033975130152874
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.039647132625252
$ws.Range("I9").Value = 1.036168659269591
$ws.Range("J9").Value = 1.036278034238513
$ws.Range("K9").Value = 1.037144327495757
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.042797858708605
$ws.Range("N9").Value = 1.016053620951186

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.029025004597021
$ws.Range("D10").Value = 1.032934365221595
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.037857653953606
$ws.Range("I10").Value = 1.035801798039948
$ws.Range("J10").Value = 1.035405380490683
$ws.Range("K10").Value = 1.036390952824614
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.041296628452782
$ws.Range("N10").Value = 1.015760755896991

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.028410818609749
$ws.Range("D11").Value = 1.03248383689641
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.037083411164268
$ws.Range("I11").Value = 1.035641601063421
$ws.Range("J11").Value = 1.035026885315719
$ws.Range("K11").Value = 1.036064024586191
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.040646502975415
$ws.Range("N11").Value = 1.015633667766922

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.028182701573332
$ws.Range("D12").Value = 1.032316511057282
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.03679591465815
$ws.Range("I12").Value = 1.035581895165898
$ws.Range("J12").Value = 1.034886201117488
$ws.Range("K12").Value = 1.035942482439907
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.040405004740572
$ws.Range("N12").Value = 1.015586420380179

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.028231632586159
$ws.Range("D13").Value = 1.032352402077638
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.03685757948964
$ws.Range("I13").Value = 1.035594711402585
$ws.Range("J13").Value = 1.034916382612576
$ws.Range("K13").Value = 1.035968558451045
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.040456807509529
$ws.Range("N13").Value = 1.015596556965833

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.028391961994638
$ws.Range("D14").Value = 1.032470005267907
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.037059644744236
$ws.Range("I14").Value = 1.035636669866535
$ws.Range("J14").Value = 1.035015258235742
$ws.Range("K14").Value = 1.036053980042768
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.040626540916118
$ws.Range("N14").Value = 1.01562976312236

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.028490748763031
$ws.Range("D15").Value = 1.032542467211655
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.03718415596101
$ws.Range("I15").Value = 1.035662495153793
$ws.Range("J15").Value = 1.035076166296956
$ws.Range("K15").Value = 1.036106597016029
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.040731117581874
$ws.Range("N15").Value = 1.015650217078761

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.029065768737861
$ws.Range("D16").Value = 1.032964268117419
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.037909050784578
$ws.Range("I16").Value = 1.035812401498789
$ws.Range("J16").Value = 1.035430486727912
$ws.Range("K16").Value = 1.036412635018854
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.041339773373421
$ws.Range("N16").Value = 1.015769184528679

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.029426496961561
$ws.Range("D17").Value = 1.033228888090423
$ws.Range("E17").Value = 0.989476357848556
$ws.Range("F17").Value = 1.038363922066092
$ws.Range("I17").Value = 1.035906074239954
$ws.Range("J17").Value = 1.035652574140169
$ws.Range("K17").Value = 1.036604414317347
$ws.Range("L17").Value = 0.9930127773699352
$ws.Range("M17").Value = 1.041721544717221
$ws.Range("N17").Value = 1.015843736024145

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.029636915529425
$ws.Range("D18").Value = 1.033383248900168
$ws.Range("E18").Value = 0.9897087662937556
$ws.Range("F18").Value = 1.038629299988447
$ws.Range("I18").Value = 1.035960582305509
$ws.Range("J18").Value = 1.035782053120192
$ws.Range("K18").Value = 1.03671620720082
$ws.Range("L18").Value = 0.9932001317071769
$ws.Range("M18").Value = 1.041944217502855
$ws.Range("N18").Value = 1.015887194041658

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.029708664828842
$ws.Range("D19").Value = 1.033435884031265
$ws.Range("E19").Value = 0.9897880325774034
$ws.Range("F19").Value = 1.038719797055504
$ws.Range("I19").Value = 1.03597914614686
$ws.Range("J19").Value = 1.035826191770978
$ws.Range("K19").Value = 1.036754314009934
$ws.Range("L19").Value = 0.9932640239640975
$ws.Range("M19").Value = 1.042020141822018
$ws.Range("N19").Value = 1.015902007582809

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.029387793003308
$ws.Range("D20").Value = 1.033200495576348
$ws.Range("E20").Value = 0.9894336180360679
$ws.Range("F20").Value = 1.038315112561897
$ws.Range("I20").Value = 1.035896037450732
$ws.Range("J20").Value = 1.035628752553258
$ws.Range("K20").Value = 1.036583845319466
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.041680585092084
$ws.Range("N20").Value = 1.015835740104911

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.028344748452494
$ws.Range("D21").Value = 1.032435373487711
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.037000139033847
$ws.Range("I21").Value = 1.035624319707957
$ws.Range("J21").Value = 1.034986144426824
$ws.Range("K21").Value = 1.036028828447126
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.040576559010431
$ws.Range("N21").Value = 1.015619985872421

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.027689054998954
$ws.Range("D22").Value = 1.031954429592371
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.036173893564367
$ws.Range("I22").Value = 1.035452313668282
$ws.Range("J22").Value = 1.034581566192055
$ws.Range("K22").Value = 1.035679251733601
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.039882340484846
$ws.Range("N22").Value = 1.015484094377498

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.028036639933616
$ws.Range("D23").Value = 1.032209375481758
$ws.Range("E23").Value = 0.9879432794643023
$ws.Range("F23").Value = 1.036611851760604
$ws.Range("I23").Value = 1.035543607790124
$ws.Range("J23").Value = 1.034796092273919
$ws.Range("K23").Value = 1.035864627135417
$ws.Range("L23").Value = 0.991776070289318
$ws.Range("M23").Value = 1.04025036590585
$ws.Range("N23").Value = 1.015556155541897

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.029405281619244
$ws.Range("D24").Value = 1.033213324892899
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.038337167292727
$ws.Range("I24").Value = 1.035900573044374
$ws.Range("J24").Value = 1.035639516691006
$ws.Range("K24").Value = 1.036593139776625
$ws.Range("L24").Value = 0.9929938892766442
$ws.Range("M24").Value = 1.041699093008076
$ws.Range("N24").Value = 1.015839353198855

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.03099397927544
$ws.Range("D25").Value = 1.034378837461393
$ws.Range("E25").Value = 0.9912096547607049
$ws.Range("F25").Value = 1.040341623565453
$ws.Range("I25").Value = 1.036309724100175
$ws.Range("J25").Value = 1.036615875015264
$ws.Range("K25").Value = 1.037435839763182
$ws.Range("L25").Value = 0.9944092447426414
$ws.Range("M25").Value = 1.043379947889332
$ws.Range("N25").Value = 1.016166943723637
